$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 0.1552006185046721
$ws.Range("I2").Value = 0.1552006185046721
$ws.Range("L2").Value = 5.871779316131055
$ws.Range("M2").Value = "[-2.1834791932672877, 13.927037825529396]"
$ws.Range("N2").Value = 0.1490177465562741
$ws.Range("O2").Value = 0.1490177465562741
$ws.Range("P2").Value = -1.308210754648002
$ws.Range("Q2").Value = "[-4.270553376951891, 1.6541318676558872]"
$ws.Range("R2").Value = 0.3784887557303631
$ws.Range("S2").Value = 0.3784887557303631
$ws.Range("T2").Value = 11.58861461336895
$ws.Range("U2").Value = "[7.147099090901431, 16.03013013583647]"
$ws.Range("V2").Value = [double]"3.916561263528351e-06"
$ws.Range("W2").Value = [double]"3.916561263528351e-06"
$ws.Range("X2").Value = 5.411331331331461
$ws.Range("Y2").Value = -6.842212212212377
$ws.Range("Z2").Value = 17.6648748748753
$ws.Range("H3").Value = 0.2677482520544663
$ws.Range("I3").Value = 0.2677482520544663
$ws.Range("L3").Value = 4.922129452800982
$ws.Range("M3").Value = "[-2.9816498715460646, 12.82590877714803]"
$ws.Range("N3").Value = 0.2162108008165877
$ws.Range("O3").Value = 0.2162108008165877
$ws.Range("P3").Value = -1.74847398938531
$ws.Range("Q3").Value = "[-4.8680534812382374, 1.3711055024676169]"
$ws.Range("R3").Value = 0.2649339031954132
$ws.Range("S3").Value = 0.2649339031954132
$ws.Range("T3").Value = 10.68822447633932
$ws.Range("U3").Value = "[6.316519089011899, 15.059929863666746]"
$ws.Range("V3").Value = [double]"1.183478876187749e-05"
$ws.Range("W3").Value = [double]"1.183478876187749e-05"
$ws.Range("X3").Value = 7.232452452452627
$ws.Range("Y3").Value = -5.671491491491624
$ws.Range("Z3").Value = 20.13639639639688
$ws.Range("H4").Value = 0.146035813440209
$ws.Range("I4").Value = 0.146035813440209
$ws.Range("L4").Value = 6.51286838247119
$ws.Range("M4").Value = "[-2.025066744415499, 15.050803509357879]"
$ws.Range("N4").Value = 0.1314451131066248
$ws.Range("O4").Value = 0.1314451131066248
$ws.Range("P4").Value = -2.64157940842385
$ws.Range("Q4").Value = "[-5.6731062533293155, 0.38994743648161556]"
$ws.Range("R4").Value = 0.08606103539380139
$ws.Range("S4").Value = 0.08606103539380139
$ws.Range("T4").Value = 13.26738160269225
$ws.Range("U4").Value = "[8.81163748673464, 17.723125718649854]"
$ws.Range("V4").Value = [double]"3.1526727872766e-07"
$ws.Range("W4").Value = [double]"3.1526727872766e-07"
$ws.Range("X4").Value = 10.92672672672699
$ws.Range("Y4").Value = -1.61299299299303
$ws.Range("Z4").Value = 23.46644644644701
$ws.Range("H5").Value = 0.6881719176464164
$ws.Range("I5").Value = 0.6881719176464164
$ws.Range("L5").Value = 2.604327285787293
$ws.Range("M5").Value = "[-5.526832454232145, 10.73548702580673]"
$ws.Range("N5").Value = 0.5221400636936955
$ws.Range("O5").Value = 0.5221400636936955
$ws.Range("P5").Value = -2.264210921506157
$ws.Range("Q5").Value = "[-5.383790413359084, 0.8553685703467702]"
$ws.Range("R5").Value = 0.1507312905353813
$ws.Range("S5").Value = 0.1507312905353813
$ws.Range("T5").Value = 10.8212161469045
$ws.Range("U5").Value = "[6.614215288960523, 15.028217004848486]"
$ws.Range("V5").Value = [double]"5.029269644740708e-06"
$ws.Range("W5").Value = [double]"5.029269644740708e-06"
$ws.Range("X5").Value = 9.365765765765989
$ws.Range("Y5").Value = -3.538178178178262
$ws.Range("Z5").Value = 22.26970970971024
$ws.Range("H6").Value = 0.5154690237693285
$ws.Range("I6").Value = 0.5154690237693285
$ws.Range("L6").Value = 3.536971533525528
$ws.Range("M6").Value = "[-4.598311946893906, 11.67225501394496]"
$ws.Range("N6").Value = 0.3858596067100342
$ws.Range("O6").Value = 0.3858596067100342
$ws.Range("P6").Value = 2.874289975356427
$ws.Range("Q6").Value = "[-0.2515789912784623, 6.000158941991317]"
$ws.Range("R6").Value = 0.07059262726495841
$ws.Range("S6").Value = 0.07059262726495841
$ws.Range("T6").Value = 10.98427494781265
$ws.Range("U6").Value = "[6.706621071864031, 15.261928823761277]"
$ws.Range("V6").Value = [double]"5.179779016595276e-06"
$ws.Range("W6").Value = [double]"5.179779016595276e-06"
$ws.Range("X6").Value = 14.10068068068102
$ws.Range("Y6").Value = 1.170720720720748
$ws.Range("Z6").Value = 27.03064064064129
$ws.Range("F7").Value = 25.14000000000049
$ws.Range("H7").Value = 0.1175585239406204
$ws.Range("I7").Value = 0.1175585239406204
$ws.Range("L7").Value = 6.164066470680401
$ws.Range("M7").Value = "[-1.6928993056220527, 14.021032246982855]"
$ws.Range("N7").Value = 0.1210802438215712
$ws.Range("O7").Value = 0.1210802438215712
$ws.Range("P7").Value = 2.232763547596349
$ws.Range("Q7").Value = "[-0.8993948938205021, 5.3649219890132]"
$ws.Range("R7").Value = 0.1579869130993461
$ws.Range("S7").Value = 0.1579869130993461
$ws.Range("T7").Value = 11.9956758034085
$ws.Range("U7").Value = "[7.788868967831366, 16.20248263898563]"
$ws.Range("V7").Value = [double]"7.499853320602767e-07"
$ws.Range("W7").Value = [double]"7.499853320602767e-07"
$ws.Range("X7").Value = 16.20636636636669
$ws.Range("Y7").Value = 3.67411411411419
$ws.Range("Z7").Value = 28.73861861861918
$ws.Range("F8").Value = 25.14000000000049
$ws.Range("H8").Value = 0.05141629019504035
$ws.Range("I8").Value = 0.05141629019504035
$ws.Range("L8").Value = 7.437821621449475
$ws.Range("M8").Value = "[-0.8076751643271756, 15.683318407226125]"
$ws.Range("N8").Value = 0.07590912531421257
$ws.Range("O8").Value = 0.07590912531421257
$ws.Range("P8").Value = 2.106974051957119
$ws.Range("Q8").Value = "[0.42139481039142357, 3.792553293522814]"
$ws.Range("R8").Value = 0.01544312383558477
$ws.Range("S8").Value = 0.01544312383558477
$ws.Range("T8").Value = 12.44923057894584
$ws.Range("U8").Value = "[8.162989280798069, 16.735471877093605]"
$ws.Range("V8").Value = [double]"5.213311780050134e-07"
$ws.Range("W8").Value = [double]"5.213311780050134e-07"
$ws.Range("X8").Value = 16.70966966966999
$ws.Range("Y8").Value = 9.965405405405594
$ws.Range("Z8").Value = 23.45393393393439
$ws.Range("F9").Value = 25.14000000000049
$ws.Range("H9").Value = 0.2236871203655357
$ws.Range("I9").Value = 0.2236871203655357
$ws.Range("L9").Value = 5.11269633362933
$ws.Range("M9").Value = "[-2.7544249618011616, 12.979817629059822]"
$ws.Range("N9").Value = 0.1972037010357215
$ws.Range("O9").Value = 0.1972037010357215
$ws.Range("P9").Value = 1.842816111114733
$ws.Range("Q9").Value = "[-1.1950002085726936, 4.8806324308021605]"
$ws.Range("R9").Value = 0.2281429672665896
$ws.Range("S9").Value = 0.2281429672665896
$ws.Range("T9").Value = 13.77425152180726
$ws.Range("U9").Value = "[9.535629551790255, 18.012873491824262]"
$ws.Range("V9").Value = [double]"4.828395505640515e-08"
$ws.Range("W9").Value = [double]"4.828395505640515e-08"
$ws.Range("X9").Value = 17.76660660660695
$ws.Range("Y9").Value = 5.611831831831942
$ws.Range("Z9").Value = 29.92138138138196
$ws.Range("F10").Value = 25.14000000000049
$ws.Range("H10").Value = 0.04013937210529694
$ws.Range("I10").Value = 0.04013937210529694
$ws.Range("L10").Value = 6.850003372350285
$ws.Range("M10").Value = "[0.12758910555652392, 13.572417639144046]"
$ws.Range("N10").Value = 0.04598137294066973
$ws.Range("O10").Value = 0.04598137294066973
$ws.Range("P10").Value = 1.817658211986887
$ws.Range("Q10").Value = "[0.40881586082750054, 3.226500563146274]"
$ws.Range("R10").Value = 0.01261195222688216
$ws.Range("S10").Value = 0.01261195222688216
$ws.Range("T10").Value = 9.924782414117061
$ws.Range("U10").Value = "[6.1214125479412065, 13.728152280292916]"
$ws.Range("V10").Value = [double]"3.908266158036966e-06"
$ws.Range("W10").Value = [double]"3.908266158036966e-06"
$ws.Range("X10").Value = 17.86726726726761
$ws.Range("Y10").Value = 12.23027027027051
$ws.Range("Z10").Value = 23.50426426426472
$ws.Range("F11").Value = 25.14000000000049
$ws.Range("H11").Value = 0.1136395342004295
$ws.Range("I11").Value = 0.1136395342004295
$ws.Range("L11").Value = 5.826031370741935
$ws.Range("M11").Value = "[-1.2397051759117712, 12.891767917395642]"
$ws.Range("N11").Value = 0.1037222134771132
$ws.Range("O11").Value = 0.1037222134771132
$ws.Range("P11").Value = 1.352237078121732
$ws.Range("Q11").Value = "[-0.3585000625718102, 3.062974218815274]"
$ws.Range("R11").Value = 0.1183791168203243
$ws.Range("S11").Value = 0.1183791168203243
$ws.Range("T11").Value = 10.98734878311639
$ws.Range("U11").Value = "[7.038262193530141, 14.936435372702643]"
$ws.Range("V11").Value = [double]"1.205144977811656e-06"
$ws.Range("W11").Value = [double]"1.205144977811656e-06"
$ws.Range("X11").Value = 19.72948948948988
$ws.Range("Y11").Value = 12.88456456456482
$ws.Range("Z11").Value = 26.57441441441494
$ws.Range("F12").Value = 25.14000000000049
$ws.Range("H12").Value = 0.0710158167728111
$ws.Range("I12").Value = 0.0710158167728111
$ws.Range("L12").Value = 7.084615070840515
$ws.Range("M12").Value = "[-1.0919869936309468, 15.261217135311977]"
$ws.Range("N12").Value = 0.08779014667207607
$ws.Range("O12").Value = 0.08779014667207607
$ws.Range("P12").Value = 1.100658086843271
$ws.Range("Q12").Value = "[-0.5975001042863468, 2.7988162779728896]"
$ws.Range("R12").Value = 0.1983791969334658
$ws.Range("S12").Value = 0.1983791969334658
$ws.Range("T12").Value = 13.04400430309345
$ws.Range("U12").Value = "[8.761721908268722, 17.32628669791818]"
$ws.Range("V12").Value = [double]"1.967440603767301e-07"
$ws.Range("W12").Value = [double]"1.967440603767301e-07"
$ws.Range("X12").Value = 20.7360960960965
$ws.Range("Y12").Value = 1.170720720720748
$ws.Range("Z12").Value = 27.53069069069123

Write-Host "Applied cosinor per-day updates"
